$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last row (Dorian Finney-Smith's row, formerly row 10) so the
# table shrinks from A1:C18 to A1:C17.
$ws.Rows.Item(18).Delete()

# New ordering / content for the player table (rows 2-17).
$players = @(
    @("Keyonte George", "PG,SG", "Utah Jazz"),
    @("Shai Gilgeous-Alexander", "PG", "Oklahoma City Thunder"),
    @("Dennis Schröder", "PG", "Brooklyn Nets"),
    @("Christian Braun", "SG,SF", "Denver Nuggets"),
    @("Brandon Boston Jr.", "SG,SF", "New Orleans Pelicans"),
    @("John Collins", "PF,C", "Utah Jazz"),
    @("Jalen Williams", "SG,SF,PF,C", "Oklahoma City Thunder"),
    @("Kyrie Irving", "PG,SG", "Dallas Mavericks"),
    @("Lauri Markkanen", "SF,PF", "Utah Jazz"),
    @("Jimmy Butler", "SF,PF", "Miami Heat"),
    @("Tobias Harris", "SF,PF", "Detroit Pistons"),
    @("Jordan Poole", "PG,SG", "Washington Wizards"),
    @("RJ Barrett", "SF,PF", "Toronto Raptors"),
    @("Zach LaVine", "SG,SF", "Chicago Bulls"),
    @("Joel Embiid", "C", "Philadelphia 76ers"),
    @("CJ McCollum", "PG,SG", "New Orleans Pelicans")
)

$row = 2
foreach ($player in $players) {
    $ws.Cells.Item($row, 1).Value = $player[0]
    $ws.Cells.Item($row, 2).Value = $player[1]
    $ws.Cells.Item($row, 3).Value = $player[2]
    $row = $row + 1
}
